$wb = $excel.ActiveWorkbook

# Add sheet "test5" at the end
$count = $wb.Worksheets.Count
$sheet5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$sheet5.Name = "test5"

$sheet5.Cells.Item(1, 1).Value = 'Word'
$sheet5.Cells.Item(1, 2).Value = 'Hypernym 1'
$sheet5.Cells.Item(1, 3).Value = 'Hypernym 2'
$sheet5.Cells.Item(1, 4).Value = 'Hypernym 3'
$sheet5.Cells.Item(1, 5).Value = 'Hypernym 4'
$sheet5.Cells.Item(1, 6).Value = 'Hypernym 5'
$sheet5.Cells.Item(2, 1).Value = 'cow'
$sheet5.Cells.Item(2, 2).Value = 'Synset(''entity.n.01'')'
$sheet5.Cells.Item(2, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet5.Cells.Item(2, 4).Value = 'Synset(''object.n.01'')'
$sheet5.Cells.Item(2, 5).Value = 'Synset(''whole.n.02'')'
$sheet5.Cells.Item(2, 6).Value = 'Synset(''living_thing.n.01'')'
$sheet5.Cells.Item(3, 1).Value = 'hello'
$sheet5.Cells.Item(3, 2).Value = 'Synset(''entity.n.01'')'
$sheet5.Cells.Item(3, 3).Value = 'Synset(''abstraction.n.06'')'
$sheet5.Cells.Item(3, 4).Value = 'Synset(''communication.n.02'')'
$sheet5.Cells.Item(3, 5).Value = 'Synset(''message.n.02'')'
$sheet5.Cells.Item(3, 6).Value = 'Synset(''acknowledgment.n.03'')'

# Add sheet "test6" at the end
$count = $wb.Worksheets.Count
$sheet6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$sheet6.Name = "test6"

$sheet6.Cells.Item(1, 1).Value = 'Word'
$sheet6.Cells.Item(1, 2).Value = 'Hypernym 1'
$sheet6.Cells.Item(1, 3).Value = 'Hypernym 2'
$sheet6.Cells.Item(1, 4).Value = 'Hypernym 3'
$sheet6.Cells.Item(1, 5).Value = 'Hypernym 4'
$sheet6.Cells.Item(1, 6).Value = 'Hypernym 5'
$sheet6.Cells.Item(2, 1).Value = 'light'
$sheet6.Cells.Item(2, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(2, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(2, 4).Value = 'Synset(''process.n.06'')'
$sheet6.Cells.Item(2, 5).Value = 'Synset(''phenomenon.n.01'')'
$sheet6.Cells.Item(2, 6).Value = 'Synset(''natural_phenomenon.n.01'')'
$sheet6.Cells.Item(3, 1).Value = 'ship'
$sheet6.Cells.Item(3, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(3, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(3, 4).Value = 'Synset(''object.n.01'')'
$sheet6.Cells.Item(3, 5).Value = 'Synset(''whole.n.02'')'
$sheet6.Cells.Item(3, 6).Value = 'Synset(''artifact.n.01'')'
$sheet6.Cells.Item(4, 1).Value = 'hello'
$sheet6.Cells.Item(4, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(4, 3).Value = 'Synset(''abstraction.n.06'')'
$sheet6.Cells.Item(4, 4).Value = 'Synset(''communication.n.02'')'
$sheet6.Cells.Item(4, 5).Value = 'Synset(''message.n.02'')'
$sheet6.Cells.Item(4, 6).Value = 'Synset(''acknowledgment.n.03'')'
$sheet6.Cells.Item(5, 1).Value = 'dfbfdbs'
$sheet6.Cells.Item(6, 1).Value = 'word'
$sheet6.Cells.Item(6, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(6, 3).Value = 'Synset(''abstraction.n.06'')'
$sheet6.Cells.Item(6, 4).Value = 'Synset(''relation.n.01'')'
$sheet6.Cells.Item(6, 5).Value = 'Synset(''part.n.01'')'
$sheet6.Cells.Item(6, 6).Value = 'Synset(''language_unit.n.01'')'
$sheet6.Cells.Item(7, 1).Value = 'walrus'
$sheet6.Cells.Item(7, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(7, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(7, 4).Value = 'Synset(''object.n.01'')'
$sheet6.Cells.Item(7, 5).Value = 'Synset(''whole.n.02'')'
$sheet6.Cells.Item(7, 6).Value = 'Synset(''living_thing.n.01'')'
$sheet6.Cells.Item(8, 1).Value = 'window'
$sheet6.Cells.Item(8, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(8, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(8, 4).Value = 'Synset(''object.n.01'')'
$sheet6.Cells.Item(8, 5).Value = 'Synset(''whole.n.02'')'
$sheet6.Cells.Item(8, 6).Value = 'Synset(''artifact.n.01'')'
$sheet6.Cells.Item(9, 1).Value = 'chair'
$sheet6.Cells.Item(9, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(9, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(9, 4).Value = 'Synset(''object.n.01'')'
$sheet6.Cells.Item(9, 5).Value = 'Synset(''whole.n.02'')'
$sheet6.Cells.Item(9, 6).Value = 'Synset(''artifact.n.01'')'
$sheet6.Cells.Item(10, 1).Value = 'photon'
$sheet6.Cells.Item(10, 2).Value = 'Synset(''entity.n.01'')'
$sheet6.Cells.Item(10, 3).Value = 'Synset(''physical_entity.n.01'')'
$sheet6.Cells.Item(10, 4).Value = 'Synset(''object.n.01'')'
$sheet6.Cells.Item(10, 5).Value = 'Synset(''whole.n.02'')'
$sheet6.Cells.Item(10, 6).Value = 'Synset(''natural_object.n.01'')'
